# Daten aktualisiert am 2023-07-01
# Appends a new block of FTSE100-style ticker symbols to Sheet1,
# extending the existing list from row 5918 through row 5980.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$tickers = @(
    "AAF",
    "ABDN",
    "ABF",
    "ANTO",
    "AUTO",
    "AV",
    "BARC",
    "BATS",
    "BDEV",
    "BEZ",
    "BF.B",
    "BKG",
    "BNZL",
    "BRBY",
    "BRK.B",
    "BT-A",
    "CCH",
    "CRDA",
    "DCC",
    "DGE",
    "ENT",
    "EXPN",
    "FCIT",
    "FRAS",
    "GLEN",
    "HLMA",
    "HSBA",
    "HSX",
    "IMB",
    "IMI",
    "INF",
    "ITRK",
    "JMAT",
    "KGF",
    "LGEN",
    "LLOY",
    "LSEG",
    "MNDI",
    "MNG",
    "OCDO",
    "PHNX",
    "PSH",
    "PSON",
    "REL",
    "RMV",
    "RR",
    "RS1",
    "SBRY",
    "SDR",
    "SGRO",
    "SKG",
    "SMDS",
    "SMT",
    "SN",
    "SPX",
    "SSE",
    "STAN",
    "STJ",
    "SVT",
    "ULVR",
    "UU",
    "WEIR",
    "WTB"
)

$startRow = 5918

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}
